$wb = $excel.ActiveWorkbook

# Update counts on both the "展览" sheet and the "全部类型" sheet.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 233
    $ws.Range("F8").Value = 7
    $ws.Range("F9").Value = 1783
    $ws.Range("F17").Value = 13933
    $ws.Range("F18").Value = 367
    $ws.Range("F22").Value = 8259
}
